# Weekly update: insert a new daily price record for "Achicoria" at
# Vega Modelo de Temuco, as row 30 (pushing the existing rows 30-59 down
# to rows 31-60, which is plain Excel row-insert semantics).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 30..59 down to 31..60, duplicating row 30's formatting into
# the freshly opened row 30.
$ws.Rows.Item(30).Insert()

# Populate the new row 30 with the new week's record.
$ws.Range("A30").Value = 10
$ws.Range("B30").Value = "Vega Modelo de Temuco"
$ws.Range("C30").Value = "La Araucanía"
$ws.Range("D30").Value = 44781
$ws.Range("E30").Value = 9
$ws.Range("F30").Value = 100112010
$ws.Range("G30").Value = "Achicoria"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 300
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 10000
$ws.Range("M30").Value = 10000
$ws.Range("N30").Value = "`$/caja 18 unidades"
$ws.Range("O30").Value = "Región Metropolitana"
$ws.Range("P30").Value = 556
$ws.Range("Q30").Value = 18
$ws.Range("R30").Value = "Hortaliza"
